$d = $word.ActiveDocument

$pairs = @(
    @("73×15=1095", "72×92=6624"),
    @("23×11=253", "69×33=2277"),
    @("28×39=1092", "34×34=1156"),
    @("60×36=2160", "72×70=5040"),
    @("67×39=2613", "99×39=3861"),
    @("36×83=2988", "75×34=2550"),
    @("93×91=8463", "32×20=640"),
    @("21×12=252", "65×70=4550"),
    @("37×58=2146", "86×69=5934"),
    @("45×33=1485", "51×29=1479"),
    @("79×97=7663", "65×49=3185"),
    @("24×44=1056", "52×69=3588"),
    @("34×76=2584", "87×35=3045"),
    @("63×71=4473", "90×49=4410"),
    @("92×38=3496", "63×42=2646"),
    @("16×17=272", "74×46=3404"),
    @("89×86=7654", "36×34=1224"),
    @("25×24=600", "27×35=945"),
    @("95×23=2185", "79×52=4108"),
    @("14×22=308", "52×70=3640"),
    @("25×61=1525", "56×82=4592"),
    @("51×32=1632", "64×15=960"),
    @("21×25=525", "47×43=2021"),
    @("26×14=364", "11×19=209"),
    @("51×80=4080", "92×23=2116")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
